# Update "想去人数" (want-to-go count) values that changed between crawl snapshots.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 83
$ws1.Range("F5").Value = 2584
$ws1.Range("F6").Value = 239

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 2

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 83
$ws4.Range("F5").Value = 2584
$ws4.Range("F6").Value = 239
$ws4.Range("F8").Value = 2
